# "kati 8 dhe 9" - add floor 8 and floor 9 entries to the building-parts
# register (Regjistri i pjeseve te nderteses), and renumber the existing
# floor-7 rows (previously mislabeled as floor 6) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) Rows 157-166 were the last block of "Kati" (floor) entries; they
#    were actually floor 7 (not floor 6 like the block above them), and
#    their "Numri i nderteses" (building part id) codes need bumping to
#    the next free sequence (...-43-0 .. ...-48-0 / ...-29-0 .. ...-32-0).
# ---------------------------------------------------------------------
$existingUpdates = @(
    @{ Row = 157; Code = "489-3-337-1-43-0"; Kati = 7 },
    @{ Row = 158; Code = "489-3-337-1-44-0"; Kati = 7 },
    @{ Row = 159; Code = "489-3-337-1-45-0"; Kati = 7 },
    @{ Row = 160; Code = "489-3-337-1-46-0"; Kati = 7 },
    @{ Row = 161; Code = "489-3-337-1-47-0"; Kati = 7 },
    @{ Row = 162; Code = "489-3-337-1-48-0"; Kati = 7 },
    @{ Row = 163; Code = "489-3-337-2-29-0"; Kati = 7 },
    @{ Row = 164; Code = "489-3-337-2-30-0"; Kati = 7 },
    @{ Row = 165; Code = "489-3-337-2-31-0"; Kati = 7 },
    @{ Row = 166; Code = "489-3-337-2-32-0"; Kati = 7 }
)

foreach ($u in $existingUpdates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.Code
    $ws.Cells.Item($u.Row, 7).Value = $u.Kati
}

# ---------------------------------------------------------------------
# 2) Append 20 new rows (167-186): ten for floor 8, ten for floor 9,
#    mirroring the apartment layout/sizes of the floor-7 block above
#    (rows 157-166) one-for-one (row 167 mirrors 157, 177 mirrors 157,
#    etc.) so each new row's "Sipërfaqja" (area, column H) is identical
#    to its source row. Copy the source row first so borders/number
#    formats/styles/area match exactly, then overwrite just the cells
#    that actually differ (Nr, Numri i nderteses code, Hyrja, Kati).
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 167; Src = 157; Nr = 160; Code = "489-3-337-1-49-0"; Hyrja = 1; Kati = 8 },
    @{ Row = 168; Src = 158; Nr = 161; Code = "489-3-337-1-50-0"; Hyrja = 1; Kati = 8 },
    @{ Row = 169; Src = 159; Nr = 162; Code = "489-3-337-1-51-0"; Hyrja = 1; Kati = 8 },
    @{ Row = 170; Src = 160; Nr = 163; Code = "489-3-337-1-52-0"; Hyrja = 1; Kati = 8 },
    @{ Row = 171; Src = 161; Nr = 164; Code = "489-3-337-1-53-0"; Hyrja = 1; Kati = 8 },
    @{ Row = 172; Src = 162; Nr = 165; Code = "489-3-337-1-54-0"; Hyrja = 1; Kati = 8 },
    @{ Row = 173; Src = 163; Nr = 166; Code = "489-3-337-2-33-0"; Hyrja = 2; Kati = 8 },
    @{ Row = 174; Src = 164; Nr = 167; Code = "489-3-337-2-34-0"; Hyrja = 2; Kati = 8 },
    @{ Row = 175; Src = 165; Nr = 168; Code = "489-3-337-2-35-0"; Hyrja = 2; Kati = 8 },
    @{ Row = 176; Src = 166; Nr = 169; Code = "489-3-337-2-36-0"; Hyrja = 2; Kati = 8 },
    @{ Row = 177; Src = 157; Nr = 170; Code = "489-3-337-1-55-0"; Hyrja = 1; Kati = 9 },
    @{ Row = 178; Src = 158; Nr = 171; Code = "489-3-337-1-56-0"; Hyrja = 1; Kati = 9 },
    @{ Row = 179; Src = 159; Nr = 172; Code = "489-3-337-1-57-0"; Hyrja = 1; Kati = 9 },
    @{ Row = 180; Src = 160; Nr = 173; Code = "489-3-337-1-58-0"; Hyrja = 1; Kati = 9 },
    @{ Row = 181; Src = 161; Nr = 174; Code = "489-3-337-1-59-0"; Hyrja = 1; Kati = 9 },
    @{ Row = 182; Src = 162; Nr = 175; Code = "489-3-337-1-60-0"; Hyrja = 1; Kati = 9 },
    @{ Row = 183; Src = 163; Nr = 176; Code = "489-3-337-2-37-0"; Hyrja = 2; Kati = 9 },
    @{ Row = 184; Src = 164; Nr = 177; Code = "489-3-337-2-38-0"; Hyrja = 2; Kati = 9 },
    @{ Row = 185; Src = 165; Nr = 178; Code = "489-3-337-2-39-0"; Hyrja = 2; Kati = 9 },
    @{ Row = 186; Src = 166; Nr = 179; Code = "489-3-337-2-40-0"; Hyrja = 2; Kati = 9 }
)

foreach ($n in $newRows) {
    $srcRange = $ws.Range("B" + $n.Src + ":Q" + $n.Src)
    $dstRange = $ws.Range("B" + $n.Row + ":Q" + $n.Row)
    $srcRange.Copy($dstRange)

    $ws.Cells.Item($n.Row, 2).Value = $n.Nr
    $ws.Cells.Item($n.Row, 3).Value = $n.Code
    $ws.Cells.Item($n.Row, 6).Value = $n.Hyrja
    $ws.Cells.Item($n.Row, 7).Value = $n.Kati
}

# ---------------------------------------------------------------------
# 3) Move the viewport/selection to the newly added block, as the
#    author left it after finishing the edit.
# ---------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("I177:Q186").Select() | Out-Null
